# Updates cryptos list values (price/volume) plus a handful of coin
# name/link/price swaps, per the Fri Nov 15 20:30:40 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.994.57'
$ws.Range("E2").Value = '  +2.78%  '

$ws.Range("D3").Value = '3.069.30'
$ws.Range("E3").Value = '  -1.74%  '

$ws.Range("D4").Value = '''1.01'
$ws.Range("E4").Value = '  +0.60%  '

$ws.Range("D5").Value = '''214.77'
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("D6").Value = '''616.36'
$ws.Range("E6").Value = '  -2.91%  '

$ws.Range("D7").Value = '''0.375'
$ws.Range("E7").Value = '  -4.04%  '

$ws.Range("D8").Value = '''0.893'
$ws.Range("E8").Value = '  +14.49%  '

$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("D10").Value = '3.070.18'
$ws.Range("E10").Value = '  -1.60%  '

$ws.Range("D11").Value = '''0.674'
$ws.Range("E11").Value = '  +18.87%  '

$ws.Range("E12").Value = '  +5.27%  '

$ws.Range("D13").Value = '''0.0000247'
$ws.Range("E13").Value = '  -1.41%  '

$ws.Range("D14").Value = '''5.38'
$ws.Range("E14").Value = '  -0.10%  '

$ws.Range("D15").Value = '90.752.76'
$ws.Range("E15").Value = '  +2.73%  '

$ws.Range("D16").Value = '''32.86'
$ws.Range("E16").Value = '  +2.02%  '

$ws.Range("D17").Value = '3.632.54'
$ws.Range("E17").Value = '  -1.66%  '

$ws.Range("D18").Value = '3.134.25'
$ws.Range("E18").Value = '  +0.24%  '

$ws.Range("D19").Value = '''3.43'
$ws.Range("E19").Value = '  +0.70%  '

$ws.Range("D20").Value = '''0.0000225'
$ws.Range("E20").Value = '  +0.94%  '

$ws.Range("D21").Value = '''13.69'
$ws.Range("E21").Value = '  +3.47%  '

$ws.Range("D22").Value = '''431.20'
$ws.Range("E22").Value = '  +2.07%  '

$ws.Range("D23").Value = '''8.41'
$ws.Range("E23").Value = '  -0.39%  '

$ws.Range("D24").Value = '''5.08'
$ws.Range("E24").Value = '  +3.56%  '

$ws.Range("D25").Value = '''5.49'
$ws.Range("E25").Value = '  +0.37%  '

$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").Value = '''11.83'
$ws.Range("E26").Value = '  +3.14%  '

$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '''83.50'
$ws.Range("E27").Value = '  +1.08%  '

$ws.Range("D28").Value = '3.206.12'
$ws.Range("E28").Value = '  -2.37%  '

$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("D30").Value = '''0.168'
$ws.Range("E30").Value = '  +7.19%  '

$ws.Range("D31").Value = '''1.06'
$ws.Range("E31").Value = '  +6.33%  '

$ws.Range("D32").Value = '''8.74'
$ws.Range("E32").Value = '  +6.66%  '

$ws.Range("D33").Value = '''3.83'
$ws.Range("E33").Value = '  -4.95%  '

$ws.Range("D34").Value = '''511.40'
$ws.Range("E34").Value = '  +1.66%  '

$ws.Range("D35").Value = '''6.87'
$ws.Range("E35").Value = '  -0.79%  '

$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").Value = '''1.83'
$ws.Range("E36").Value = '  +0.07%  '

$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").Value = '''23.02'
$ws.Range("E37").Value = '  +2.93%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.137'
$ws.Range("E38").Value = '  -7.88%  '

$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = '''1.26'
$ws.Range("E39").Value = '  -0.89%  '

$ws.Range("D40").Value = '''22.30'
$ws.Range("E40").Value = '  +0.37%  '

$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("E43").Value = '  +4.84%  '

$ws.Range("D44").Value = '''0.367'
$ws.Range("E44").Value = '  +0.26%  '

$ws.Range("D45").Value = '''1.86'
$ws.Range("E45").Value = '  +0.50%  '

$ws.Range("D46").Value = '''0.0714'
$ws.Range("E46").Value = '  +9.27%  '

$ws.Range("D47").Value = '''144.50'
$ws.Range("E47").Value = '  -0.85%  '

$ws.Range("D48").Value = '''43.61'
$ws.Range("E48").Value = '  +0.08%  '

$ws.Range("E49").Value = '  +7.67%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '''0.000258'
$ws.Range("E50").Value = '  +8.54%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '''163.69'
$ws.Range("E51").Value = '  +0.51%  '
